# Updates the "Price" (D) and "Volume(1h)" (E) columns of the cryptos
# sheet with freshly scraped values. D-column values that look numeric
# are written as text (NumberFormat "@") so Excel does not coerce them
# into floating point numbers -- that would corrupt values which rely on
# exact decimal formatting (e.g. "88.80" staying "88.80" instead of
# becoming 88.8, or "0.000008657" staying decimal instead of becoming
# scientific notation). ClearFormats() afterwards removes the temporary
# text-format style so no stray cell formatting is left behind.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.986.74'
$ws.Range('E2').Value = '  +1.29%  '
$ws.Range('D3').Value = '1.848.66'
$ws.Range('E3').Value = '  +1.20%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.014'
$ws.Range('D4').ClearFormats()
$ws.Range('E4').Value = '  +0.44%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '309.61'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +0.26%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4777'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  +1.67%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3678'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  +2.20%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07223'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  +1.14%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.9294'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  +2.55%  '
$ws.Range('E11').Value = '  +1.57%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.07745'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  +0.97%  '
$ws.Range('D13').Value = '1.816.50'
$ws.Range('E13').Value = '  -0.25%  '
$ws.Range('E14').Value = '  +1.59%  '
$ws.Range('E15').Value = '  +1.33%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '88.80'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  +1.36%  '
$ws.Range('E17').Value = '  +0.42%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000008657'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  +0.93%  '
$ws.Range('D20').Value = '27.009.44'
$ws.Range('E20').Value = '  +1.35%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '14.46'
$ws.Range('D21').ClearFormats()
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.062'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +0.88%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '10.64'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +0.99%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.922'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +0.70%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '152.87'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +0.10%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '18.24'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +1.91%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.005'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +0.42%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '114.22'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  +0.55%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '4.955'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  +1.69%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.08876'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  +0.60%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.322'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  +5.26%  '
$ws.Range('E32').Value = '  +1.30%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.7421'
$ws.Range('D33').ClearFormats()
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.506'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  +1.79%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.745'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  -3.63%  '
$ws.Range('E36').Value = '  +3.54%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.01965'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  +1.99%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.05270'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  +2.39%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.981'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  +1.18%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.5212'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  +3.09%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '6.992'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  +2.10%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.1514'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  +1.19%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '8.244'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  +2.17%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '10.64'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +6.59%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.4741'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  +2.08%  '
$ws.Range('E46').Value = '  +0.39%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '101.92'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  +3.77%  '
$ws.Range('E48').Value = '  +2.40%  '
$ws.Range('E49').Value = '  +3.24%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.06069'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +0.62%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.8879'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  +4.14%  '
